$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Ready for handoff" -> "Handback transform failed" everywhere it appears
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# New "Error Detail" messages describing the handback filename mismatch
$wsZhCn.Range("K3").Value = "Handback file name: ztevtqrb.mto is different with handoff file name: 99bf7219-d081-4b5d-a609-40b3eacdc05b.4aeb91bec1139477220b6efb22bf5675229a5fc6.zh-cn."
$wsDeDe.Range("K3").Value = "Handback file name: ztevtqrb.mto is different with handoff file name: 99bf7219-d081-4b5d-a609-40b3eacdc05b.4aeb91bec1139477220b6efb22bf5675229a5fc6.de-de."
